$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values are stored as plain text in the sheet (e.g. "25.914.92",
# "0.06390", multi-dot thousand groupings, unicode subscript digits, etc.). A leading
# apostrophe forces Excel to keep them as literal text instead of re-parsing them as
# numbers (which would silently drop trailing zeros / mangle the multi-dot values).

$ws.Range("D2").Value = "'25.914.92"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "'1.640.57"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -1.23%  "

$ws.Range("D5").Value = "'215.24"
$ws.Range("E5").Value = "  -0.67%  "

$ws.Range("D6").Value = "'0.5038"
$ws.Range("E6").Value = "  -0.45%  "

$ws.Range("E7").Value = "  -1.10%  "

$ws.Range("D8").Value = "'0.2571"
$ws.Range("E8").Value = "  -0.89%  "

$ws.Range("D9").Value = "'0.06390"
$ws.Range("E9").Value = "  -1.17%  "

$ws.Range("E10").Value = "  +0.18%  "

$ws.Range("D11").Value = "'0.07780"
$ws.Range("E11").Value = "  +0.11%  "

$ws.Range("D12").Value = "'1.655.21"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "'4.278"
$ws.Range("E13").Value = "  -0.09%  "

$ws.Range("D14").Value = "'1.864.09"
$ws.Range("E14").Value = "  -0.65%  "

$ws.Range("D15").Value = "'0.5434"
$ws.Range("E15").Value = "  -1.17%  "

$ws.Range("D16").Value = "'0.0₅7858"
$ws.Range("E16").Value = "  -1.48%  "

$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "'25.950.16"
$ws.Range("E18").Value = "  -0.78%  "

$ws.Range("E19").Value = "  -0.89%  "

$ws.Range("D20").Value = "'198.84"
$ws.Range("E20").Value = "  -2.75%  "

$ws.Range("D21").Value = "'4.389"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("D22").Value = "'9.948"
$ws.Range("E22").Value = "  -1.03%  "

$ws.Range("D23").Value = "'5.977"
$ws.Range("E23").Value = "  -0.40%  "

$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("D25").Value = "'1.876"
$ws.Range("E25").Value = "  -4.40%  "

$ws.Range("D26").Value = "'140.03"
$ws.Range("E26").Value = "  -1.88%  "

$ws.Range("D27").Value = "'0.1144"
$ws.Range("E27").Value = "  -1.66%  "

$ws.Range("D28").Value = "'6.859"
$ws.Range("E28").Value = "  +1.15%  "

$ws.Range("D29").Value = "'15.71"
$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("D30").Value = "'1.243"
$ws.Range("E30").Value = "  -0.54%  "

$ws.Range("D31").Value = "'0.04895"
$ws.Range("E31").Value = "  -3.88%  "

$ws.Range("D32").Value = "'3.258"
$ws.Range("E32").Value = "  -0.85%  "

$ws.Range("D33").Value = "'3.193"
$ws.Range("E33").Value = "  -0.74%  "

$ws.Range("D34").Value = "'1.531"
$ws.Range("E34").Value = "  -1.64%  "

$ws.Range("D35").Value = "'2.370"
$ws.Range("E35").Value = "  +0.65%  "

$ws.Range("D36").Value = "'0.8928"
$ws.Range("E36").Value = "  -1.08%  "

$ws.Range("D37").Value = "'2.604"
$ws.Range("E37").Value = "  -1.33%  "

$ws.Range("D38").Value = "'1.139.58"
$ws.Range("E38").Value = "  -1.66%  "

$ws.Range("D39").Value = "'0.5549"
$ws.Range("E39").Value = "  -2.20%  "

$ws.Range("D40").Value = "'0.01560"
$ws.Range("E40").Value = "  -1.35%  "

$ws.Range("E41").Value = "  -0.92%  "

$ws.Range("D42").Value = "'5.689"
$ws.Range("E42").Value = "  -0.10%  "

$ws.Range("D43").Value = "'0.8170"
$ws.Range("E43").Value = "  -0.53%  "

$ws.Range("D44").Value = "'99.59"
$ws.Range("E44").Value = "  -0.77%  "

$ws.Range("D45").Value = "'0.0₈119"
$ws.Range("E45").Value = "  +4.73%  "

$ws.Range("D46").Value = "'1.774.80"
$ws.Range("E46").Value = "  -0.59%  "

$ws.Range("D47").Value = "'0.4518"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("D48").Value = "'1.009"
$ws.Range("E48").Value = "  -0.46%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").Value = "'0.05085"
$ws.Range("E50").Value = "  +0.57%  "

$ws.Range("E51").Value = "  -0.68%  "
